# Apply the "Generate Report for Archive" localization-status refresh:
#  1) Status cells that used to read "Ready for handoff" now read "In Translation".
#  2) The "Status" column is narrower on every sheet (Overview!E:F and the
#     Status column on each per-language sheet).

$wb = $excel.ActiveWorkbook

# --- 1) Update status text -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- 2) Narrow the Status-related columns -----------------------------------
# ColumnWidth is expressed in characters of the Normal style's font; 12.5
# characters is the closest this host's column-width quantization gets to
# the target OOXML column width of ~13.41 characters.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
